# Case_5_52 (380 kV case): update res_line/pl_mw.xlsx values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.279501564648825
$ws.Range("C2").Value = 0.03613163393042385
$ws.Range("D2").Value = 0.007396171535118512
$ws.Range("E2").Value = 0.07614127728833964
$ws.Range("F2").Value = 4.461247904343765
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1909334364561239
$ws.Range("K2").Value = 1.709876221726944
$ws.Range("L2").Value = 0.3211692963923198
$ws.Range("N2").Value = 4.371269726001032

$ws.Range("B3").Value = 2.243400116295135
$ws.Range("C3").Value = 0.03193831572107797
$ws.Range("D3").Value = 0.007351124521136398
$ws.Range("E3").Value = 0.07627828469592934
$ws.Range("F3").Value = 4.449932950733128
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1913817980919532
$ws.Range("K3").Value = 1.672240611952446
$ws.Range("L3").Value = 0.3189368142782669
$ws.Range("N3").Value = 4.379460405894605

$ws.Range("B4").Value = 2.222451037213801
$ws.Range("C4").Value = 0.02937647600124649
$ws.Range("D4").Value = 0.007328578973478272
$ws.Range("E4").Value = 0.07638527183219423
$ws.Range("F4").Value = 4.444562865237543
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1917044313801846
$ws.Range("K4").Value = 1.650109606069549
$ws.Range("L4").Value = 0.3177171516662725
$ws.Range("N4").Value = 4.385348231184537

$ws.Range("B5").Value = 2.214220414131688
$ws.Range("C5").Value = 0.02833568208197335
$ws.Range("D5").Value = 0.007320688106432272
$ws.Range("E5").Value = 0.07643463369519132
$ws.Range("F5").Value = 4.442771181423694
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1918478270541257
$ws.Range("K5").Value = 1.641336830198725
$ws.Range("L5").Value = 0.3172581729289803
$ws.Range("N5").Value = 4.387963406231279

$ws.Range("B6").Value = 2.212872230637402
$ws.Range("C6").Value = 0.02816304913454815
$ws.Range("D6").Value = 0.007319456466948537
$ws.Range("E6").Value = 0.07644317873528905
$ws.Range("F6").Value = 4.442497632788857
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1918723582268989
$ws.Range("K6").Value = 1.63989496305166
$ws.Range("L6").Value = 0.3171842591205873
$ws.Range("N6").Value = 4.388410687540528

$ws.Range("B7").Value = 2.222338795725364
$ws.Range("C7").Value = 0.02936242670821798
$ws.Range("D7").Value = 0.007328467289729801
$ws.Range("E7").Value = 0.07638591418704266
$ws.Range("F7").Value = 4.4445370958328
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1917063169785038
$ws.Range("K7").Value = 1.649990298139869
$ws.Range("L7").Value = 0.3177108076131532
$ws.Range("N7").Value = 4.385382626573403

$ws.Range("B8").Value = 2.266801212470341
$ws.Range("C8").Value = 0.03468306769246965
$ws.Range("D8").Value = 0.007379584404757722
$ws.Range("E8").Value = 0.07618378052860564
$ws.Range("F8").Value = 4.457019137313253
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1910782164493519
$ws.Range("K8").Value = 1.696696589434652
$ws.Range("L8").Value = 0.3203682064279221
$ws.Range("N8").Value = 4.373915587436187

$ws.Range("B9").Value = 2.363650991200359
$ws.Range("C9").Value = 0.04522226193438428
$ws.Range("D9").Value = 0.007519946974543501
$ws.Range("E9").Value = 0.075968221618961
$ws.Range("F9").Value = 4.494016150565898
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.1902214516060567
$ws.Range("K9").Value = 1.796049890465724
$ws.Range("L9").Value = 0.3267766505349385
$ws.Range("N9").Value = 4.358247702251745

$ws.Range("B10").Value = 2.440706717701914
$ws.Range("C10").Value = 0.05303499218545937
$ws.Range("D10").Value = 0.00764695789959724
$ws.Range("E10").Value = 0.07591935155317131
$ws.Range("F10").Value = 4.528844567540176
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.1898197669933914
$ws.Range("K10").Value = 1.873797700398057
$ws.Range("L10").Value = 0.3322137852194516
$ws.Range("N10").Value = 4.350902866766006

$ws.Range("B11").Value = 2.477045965152286
$ws.Range("C11").Value = 0.05660550998992164
$ws.Range("D11").Value = 0.007709808263019369
$ws.Range("E11").Value = 0.07592074678250249
$ws.Range("F11").Value = 4.546353116990787
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.189686333797205
$ws.Range("K11").Value = 1.91020453904477
$ws.Range("L11").Value = 0.3348454035266286
$ws.Range("N11").Value = 4.348468425253174

$ws.Range("B12").Value = 2.490991708428623
$ws.Range("C12").Value = 0.05796002065480366
$ws.Range("D12").Value = 0.007734327528474694
$ws.Range("E12").Value = 0.07592465947683102
$ws.Range("F12").Value = 4.553222735057517
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.1896428808468222
$ws.Range("K12").Value = 1.924140485757448
$ws.Range("L12").Value = 0.3358646537959942
$ws.Range("N12").Value = 4.347677124835812

$ws.Range("B13").Value = 2.487980022870943
$ws.Range("C13").Value = 0.05766819329566886
$ws.Range("D13").Value = 0.007729015040830234
$ws.Range("E13").Value = 0.07592366649547522
$ws.Range("F13").Value = 4.551732585971024
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.1896519247501587
$ws.Range("K13").Value = 1.921132481575341
$ws.Range("L13").Value = 0.3356441304640754
$ws.Range("N13").Value = 4.347841735509846

$ws.Range("B14").Value = 2.478189585513746
$ws.Range("C14").Value = 0.05671689726024454
$ws.Range("D14").Value = 0.007711811121950518
$ws.Range("E14").Value = 0.07592100092915821
$ws.Range("F14").Value = 4.546913483877063
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.1896826171724157
$ws.Range("K14").Value = 1.911348061824015
$ws.Range("L14").Value = 0.3349288028768598
$ws.Range("N14").Value = 4.348400706639993

$ws.Range("B15").Value = 2.472216731333504
$ws.Range("C15").Value = 0.05613451976846306
$ws.Range("D15").Value = 0.007701366584882408
$ws.Range("E15").Value = 0.07591980854878955
$ws.Range("F15").Value = 4.543992839021371
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.1897023381907594
$ws.Range("K15").Value = 1.905374289296816
$ws.Range("L15").Value = 0.3344936008502799
$ws.Range("N15").Value = 4.34876010201539

$ws.Range("B16").Value = 2.438357735780585
$ws.Range("C16").Value = 0.05280198736640784
$ws.Range("D16").Value = 0.007642951674572629
$ws.Range("E16").Value = 0.07591973467317636
$ws.Range("F16").Value = 4.527733853227829
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.1898294779704592
$ws.Range("K16").Value = 1.871439339904384
$ws.Range("L16").Value = 0.3320449834096451
$ws.Range("N16").Value = 4.351080224747321

$ws.Range("B17").Value = 2.417915696213413
$ws.Range("C17").Value = 0.05076185052472226
$ws.Range("D17").Value = 0.007608407872730538
$ws.Range("E17").Value = 0.07592573129105951
$ws.Range("F17").Value = 4.518185985095641
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.1899200902538247
$ws.Range("K17").Value = 1.850887495707866
$ws.Range("L17").Value = 0.3305833350943459
$ws.Range("N17").Value = 4.352735913623093

$ws.Range("B18").Value = 2.406279045166855
$ws.Range("C18").Value = 0.04958996385984449
$ws.Range("D18").Value = 0.007589017237016549
$ws.Range("E18").Value = 0.07593140499766449
$ws.Range("F18").Value = 4.5128510084878
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1899768487155313
$ws.Range("K18").Value = 1.839164406877757
$ws.Range("L18").Value = 0.3297575301809985
$ws.Range("N18").Value = 4.353773559748959

$ws.Range("B19").Value = 2.402359871960584
$ws.Range("C19").Value = 0.04919344674075887
$ws.Range("D19").Value = 0.00758253431932765
$ws.Range("E19").Value = 0.07593370849758507
$ws.Range("F19").Value = 4.511071588463111
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1899968636037421
$ws.Range("K19").Value = 1.835211963855073
$ws.Range("L19").Value = 0.3294804865868457
$ws.Range("N19").Value = 4.354139540347944

$ws.Range("B20").Value = 2.420079256747385
$ws.Range("C20").Value = 0.05097886608966462
$ws.Range("D20").Value = 0.00761203571876834
$ws.Range("E20").Value = 0.07592486279746957
$ws.Range("F20").Value = 4.51918615301318
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.1899099642275957
$ws.Range("K20").Value = 1.853065153457067
$ws.Range("L20").Value = 0.3307373886186014
$ws.Range("N20").Value = 4.352550829361434

$ws.Range("B21").Value = 2.481060257990578
$ws.Range("C21").Value = 0.0569962494112275
$ws.Range("D21").Value = 0.007716844891804797
$ws.Range("E21").Value = 0.07592169212415989
$ws.Range("F21").Value = 4.54832246992342
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.1896734101469661
$ws.Range("K21").Value = 1.914217926031881
$ws.Range("L21").Value = 0.3351382958122571
$ws.Range("N21").Value = 4.348232978194062

$ws.Range("B22").Value = 2.521992329235729
$ws.Range("C22").Value = 0.06094316962493451
$ws.Range("D22").Value = 0.007789531154255513
$ws.Range("E22").Value = 0.07593934050798978
$ws.Range("F22").Value = 4.568760838496132
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.1895600426869208
$ws.Range("K22").Value = 1.955056042792421
$ws.Range("L22").Value = 0.3381469166725282
$ws.Range("N22").Value = 4.346172159741002

$ws.Range("B23").Value = 2.500047553940306
$ws.Range("C23").Value = 0.05883530194873288
$ws.Range("D23").Value = 0.007750357311739009
$ws.Range("E23").Value = 0.07592812109266589
$ws.Range("F23").Value = 4.557724719046035
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.1896167804160065
$ws.Range("K23").Value = 1.933180228689594
$ws.Range("L23").Value = 0.3365290601553426
$ws.Range("N23").Value = 4.347202352062666

$ws.Range("B24").Value = 2.419100750477924
$ws.Range("C24").Value = 0.05088075019720861
$ws.Range("D24").Value = 0.007610394108116481
$ws.Range("E24").Value = 0.07592524850841542
$ws.Range("F24").Value = 4.518733496718653
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.1899145276748868
$ws.Range("K24").Value = 1.852080346455978
$ws.Range("L24").Value = 0.3306676958022621
$ws.Range("N24").Value = 4.352634238794479

$ws.Range("B25").Value = 2.336415308858193
$ws.Range("C25").Value = 0.04235929651451897
$ws.Range("D25").Value = 0.007477739695907815
$ws.Range("E25").Value = 0.07600725391176155
$ws.Range("F25").Value = 4.482665705421567
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1904131722146971
$ws.Range("K25").Value = 1.768339057237341
$ws.Range("L25").Value = 0.3249149307799257
$ws.Range("N25").Value = 4.361755317418897
